# Fix table report for requirements
# Re-orders the functional-requirements rows (B2:B4) and the non-functional
# requirements table (A20:C29) to match the corrected list, and grows the
# "Tabella2" table by one row to hold the new NFR10 entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Functional requirements (table "Tabella1", A1:C17) ---
# Row 2 used to hold a stray requirement that really belongs with the NFRs;
# shift R1/R2 up one slot and add the new "redirect when not logged in"
# requirement as R3.
$ws.Range("B2").Value = "User shall be able to register in the Web App"
$ws.Range("B3").Value = "User shall be able to login in the Web App"
$ws.Range("B4").Value = "User should be redirected to homepage if not logged in"

# --- Non-functional requirements (table "Tabella2", A19:C28 -> A19:C29) ---
# Insert the "display home page for unregistered users" requirement at the
# top (NFR1) and push the rest down by one row, appending NFR10 at the end.
$ws.Range("B20").Value = "Web App shall display the home page for users not registered"
$ws.Range("B21").Value = "Web App shall have data integrity meaning that the data is stored persistently"
$ws.Range("B22").Value = "Web App shall store different details associated to different user"
$ws.Range("B23").Value = "Web App shall be accessible from any smartphone"
$ws.Range("B24").Value = "Web App shall be working with Android NFC"
$ws.Range("C24").Value = "MUST have"
$ws.Range("B25").Value = "Web App shall be able to show the user a collection of the stamp in the UI"
$ws.Range("C25").NumberFormat = "0.00"
$ws.Range("C25").Value = "SHOULD have"
$ws.Range("B26").Value = "Web App shall be visually appealing to the customer"
$ws.Range("C26").NumberFormat = "0.00"
$ws.Range("C26").Value = "COULD have"
$ws.Range("B27").Value = "Web App shall render well on mobile devices"
$ws.Range("C27").Value = "COULD have"
$ws.Range("B28").Value = "Web App shall have different login depending on customer or cashier/retailer"

# New row 29 holding the final NFR (NFR10)
$ws.Range("A29").Value = "NFR10"
$ws.Range("B29").Value = "Web App shall have different content depending on customer or cashier/retailer"
$ws.Range("C29").Value = "WON'T have"

# Grow the NFR table so the new row is included
$lo2 = $ws.ListObjects.Item("Tabella2")
$lo2.Resize($ws.Range("A19:C29"))

# Visually mark the new top row of the NFR table with top/bottom borders
$topBorder = $ws.Range("B20").Borders.Item(8)
$topBorder.LineStyle = 1
$topBorder.Weight = 2
$bottomBorder = $ws.Range("B20").Borders.Item(9)
$bottomBorder.LineStyle = 1
$bottomBorder.Weight = 2

# Column width adjustments (content changed length so the best-fit columns
# were resized)
$ws.Columns.Item(1).ColumnWidth = 5.666666666666667
$ws.Columns.Item(2).ColumnWidth = 67.66666666666667

# Restore the view: no frozen scroll position, new zoom level and selection
$ws.Activate()
$excel.ActiveWindow.Zoom = 145
$ws.Range("E4").Select()
